$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new numeric value next to the "Bold" label, using the bold font
# and a "0.00" number format (new cellXf #7 in styles.xml).
$ws.Range("B2").Value = 1.26
$ws.Range("B2").NumberFormat = "0.00"

# Leave the selection on the newly entered cell.
$ws.Range("B2").Select() | Out-Null
